# Weekly price update: insert a new observation row for "Perejil" (Femacal de
# La Calera) at row 18, pushing the previously existing rows 18-22 down to
# rows 19-23 (the data itself is unchanged for those rows, only their row
# number shifts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; existing rows 18-22 shift to 19-23.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with this week's data.
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "Femacal de La Calera"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44559
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 100112044
$ws.Range("G18").Value = "Perejil"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 2000
$ws.Range("N18").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O18").Value = "Provincia de Quillota"
$ws.Range("P18").Value = 667
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = "Hortaliza"
